$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the account 004426743 / GABRIELLE / 10000 row entirely ---
$ws.Rows.Item(14).Delete()

# --- Insert new row for account 005003629 / ANDRE / 1998.87 ---
# (keeps the Saldo column sorted descending; goes just above
#  005186167 / ANDREA / 1674.77, now at row 20 after the delete above)
$ws.Rows.Item(20).Insert()
$ws.Cells.Item(20, 1).Value = "'005003629"
$ws.Cells.Item(20, 2).Value = "ANDRE"
$ws.Cells.Item(20, 3).Value = 1998.87

# --- Remove the stale account 004228456 / FLASH / 73.98 row ---
# (still row 96 -- the earlier delete/insert pair above cancel out in
#  row count before this point)
$ws.Rows.Item(96).Delete()

# --- Insert the updated account 004228456 / FLASH / 629.61 row ---
# goes just above 005079311 / JOVINO / 623.67, which is row 32 at this point
$ws.Rows.Item(32).Insert()
$ws.Cells.Item(32, 1).Value = "'004228456"
$ws.Cells.Item(32, 2).Value = "FLASH"
$ws.Cells.Item(32, 3).Value = 629.61
